$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Sergipe", "Taxa de estupro", "31/12/2012", 22.0288630216873),
    @("Sergipe", "Taxa de estupro", "31/12/2013", 22.68108661533515),
    @("Sergipe", "Taxa de estupro", "31/12/2014", 22.97738214630375),
    @("Sergipe", "Taxa de estupro", "31/12/2015", 19.75088912439359),
    @("Sergipe", "Taxa de estupro", "31/12/2016", 16.46232929160346),
    @("Sergipe", "Taxa de estupro", "31/12/2017", 24.82391626997932),
    @("Sergipe", "Taxa de estupro", "31/12/2018", 23.96515308729109),
    @("Sergipe", "Taxa de estupro", "31/12/2019", 31.84414120005429),
    @("Sergipe", "Taxa de estupro", "31/12/2020", 26.60833819930982),
    @("Sergipe", "Taxa de estupro", "31/12/2021", 33.73995178052012)
)

$row = 12
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row = $row + 1
}
